# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Mango, Brasil) right before the current
# row 217, pushing the existing rows 217-263 down to 219-265.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 217-218 (shifts old 217..263 -> 219..265)
$ws.Range("A217:A218").EntireRow.Insert()

# New row 217 data
$row217 = @(
    4,
    "Feria Lagunitas de Puerto Montt",
    "Los Lagos",
    44855,
    10,
    "Fruta",
    100108,
    "Tropicales y subtropicales",
    100108002,
    "Mango",
    "Sin especificar",
    "Primera",
    120,
    9000,
    10000,
    9500,
    "$/bandeja 4 kilos",
    "Brasil",
    2375,
    4
)

# New row 218 data
$row218 = @(
    4,
    "Feria Lagunitas de Puerto Montt",
    "Los Lagos",
    44855,
    10,
    "Fruta",
    100108,
    "Tropicales y subtropicales",
    100108002,
    "Mango",
    "Sin especificar",
    "Segunda",
    100,
    8000,
    8000,
    8000,
    "$/bandeja 4 kilos",
    "Brasil",
    2000,
    4
)

for ($i = 0; $i -lt $row217.Length; $i++) {
    $ws.Cells.Item(217, $i + 1).Value2 = $row217[$i]
}

for ($i = 0; $i -lt $row218.Length; $i++) {
    $ws.Cells.Item(218, $i + 1).Value2 = $row218[$i]
}
